$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 203, pushing the
# existing rows 203-214 down to 205-216.
$ws.Rows.Item(203).Resize(2).Insert()

# Row 203: new Red Globe record (Region de O'Higgins)
$ws.Cells.Item(203, 1).Value = 11
$ws.Cells.Item(203, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(203, 3).Value = "Bíobío"
$ws.Cells.Item(203, 4).Value = 45021
$ws.Cells.Item(203, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(203, 5).Value = 8
$ws.Cells.Item(203, 6).Value = "Fruta"
$ws.Cells.Item(203, 7).Value = 100109
$ws.Cells.Item(203, 8).Value = "Uva"
$ws.Cells.Item(203, 9).Value = 100109001
$ws.Cells.Item(203, 10).Value = "Uva"
$ws.Cells.Item(203, 11).Value = "Red Globe"
$ws.Cells.Item(203, 12).Value = "Primera"
$ws.Cells.Item(203, 13).Value = 220
$ws.Cells.Item(203, 14).Value = 11000
$ws.Cells.Item(203, 15).Value = 12000
$ws.Cells.Item(203, 16).Value = 11545
$ws.Cells.Item(203, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(203, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(203, 19).Value = 641
$ws.Cells.Item(203, 20).Value = 18

# Row 204: new Thompson seedless record (Region de O'Higgins)
$ws.Cells.Item(204, 1).Value = 11
$ws.Cells.Item(204, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(204, 3).Value = "Bíobío"
$ws.Cells.Item(204, 4).Value = 45021
$ws.Cells.Item(204, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(204, 5).Value = 8
$ws.Cells.Item(204, 6).Value = "Fruta"
$ws.Cells.Item(204, 7).Value = 100109
$ws.Cells.Item(204, 8).Value = "Uva"
$ws.Cells.Item(204, 9).Value = 100109001
$ws.Cells.Item(204, 10).Value = "Uva"
$ws.Cells.Item(204, 11).Value = "Thompson seedless"
$ws.Cells.Item(204, 12).Value = "Primera"
$ws.Cells.Item(204, 13).Value = 150
$ws.Cells.Item(204, 14).Value = 11000
$ws.Cells.Item(204, 15).Value = 12000
$ws.Cells.Item(204, 16).Value = 11467
$ws.Cells.Item(204, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(204, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(204, 19).Value = 637
$ws.Cells.Item(204, 20).Value = 18
